{"js": "// READ ME.docx update \u2014 Office.js (Word JavaScript API) version.\n//\n// Changes applied (see commit \"Updating the READ ME files\"):\n//  1. Title paragraph text is consolidated (\"Three\"/\" IBR\"/\"s\" runs merged) \u2014\n//     no visible text change, just normalizing the run to a single run.\n//  2. New sub-bullet inserted after \"PSSE version 35 used for results shown\"\n//     explaining the `eventType` switch in disturbance.py (PSSE section).\n//  3. \"EPCL (.p) files ...\" bullet consolidated into a single run (no visible\n//     text change).\n//  4. New sub-bullet inserted after \"PSLF version 32 used for results shown\"\n//     explaining the `case_idx` switch in Disturbance.p (PSLF section).\n//  5. New sub-bullet inserted after \"PSCAD version 5.0 used for results shown\"\n//     (before the E-TRAN bullet) explaining the Timed Breaker Logic switch,\n//     and the E-TRAN bullet's two runs (\"...PSCAD fil\" + \"e\") are\n//     consolidated into one run (no visible text change).\n\nconst LDQ = \"\\u201C\"; // \u201c\nconst RDQ = \"\\u201D\"; // \u201d\nconst LSQ = \"\\u2018\"; // \u2018\nconst RSQ = \"\\u2019\"; // \u2019\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Helper: find the (first) paragraph whose text matches exactly.\nfunction findParagraph(text) {\n  const p = paragraphs.items.find((pp) => pp.text === text);\n  if (!p) {\n    throw new Error(\"Could not locate paragraph with text: \" + text);\n  }\n  return p;\n}\n\n// --- 1. Normalize the title paragraph to a single run ------------------\nconst title = paragraphs.items[0];\ntitle.load(\"text\");\nawait context.sync();\nif (title.text === \"READ ME: Enhanced IEEE 39 Bus System with Three IBRs\") {\n  title.clear();\n  title.insertText(\"READ ME: Enhanced IEEE 39 Bus System with Three IBRs\", \"Start\");\n}\n\n// --- 2. Insert the eventType/disturbance.py bullet (PSSE section) ------\nconst psseVersionPara = findParagraph(\"PSSE version 35 used for results shown \");\nconst eventTypePara = psseVersionPara.insertParagraph(\n  \"By modifying the value of \" +\n    LDQ + \"eventType\" + RDQ +\n    \" in the \" + LSQ + \"disturbance.py\" + RSQ +\n    \" file, one can simulate either a Busfault or GenTrip scenario\",\n  \"After\"\n);\neventTypePara.listItem.level = 1;\n\n// --- 3. Normalize the EPCL (.p) bullet to a single run ------------------\nconst epclPara = findParagraph(\"EPCL (.p) files with for running dynamic simulation with run commands\");\nepclPara.clear();\nepclPara.insertText(\"EPCL (.p) files with for running dynamic simulation with run commands\", \"Start\");\n\n// --- 4. Insert the case_idx/Disturbance.p bullet (PSLF section) --------\nconst pslfVersionPara = findParagraph(\"PSLF version 32 used for results shown\");\nconst caseIdxPara = pslfVersionPara.insertParagraph(\n  \"By modifying the value of \" +\n    LDQ + \"case_idx\" + RDQ +\n    \" in the \" + LSQ + \"Disturbance.p\" + RSQ +\n    \" file, one can simulate either a BusFault or GenTrip scenario\",\n  \"After\"\n);\ncaseIdxPara.listItem.level = 1;\n\n// --- 5. Insert the Timed Breaker Logic bullet (PSCAD section) + fix ETRAN bullet\nconst pscadVersionPara = findParagraph(\"PSCAD version 5.0 used for results shown\");\nconst breakerPara = pscadVersionPara.insertParagraph(\n  \"By modifying the time of breaker operation in \" +\n    LDQ + \"Timed Breaker Logic\" + RDQ +\n    \" in PSCAD, one can simulate either a BusFault or GenTrip scenario\",\n  \"After\"\n);\nbreakerPara.listItem.level = 1;\n\nconst etranPara = findParagraph(\n  \"E-TRAN runtime library file for initializing ETRAN-based components in PSCAD file\"\n);\netranPara.clear();\netranPara.insertText(\n  \"E-TRAN runtime library file for initializing ETRAN-based components in PSCAD file\",\n  \"Start\"\n);\n\nawait context.sync();\n", "ps1": "# READ ME.docx update - Word COM interop (PowerShell-style) version.\n#\n# Changes applied (see commit \"Updating the READ ME files\"):\n#  1. Title paragraph text is consolidated (\"Three\"/\" IBR\"/\"s\" runs merged) -\n#     no visible text change, just normalizing the run to a single run.\n#  2. New sub-bullet inserted after \"PSSE version 35 used for results shown\"\n#     explaining the `eventType` switch in disturbance.py (PSSE section).\n#  3. \"EPCL (.p) files ...\" bullet consolidated into a single run (no visible\n#     text change).\n#  4. New sub-bullet inserted after \"PSLF version 32 used for results shown\"\n#     explaining the `case_idx` switch in Disturbance.p (PSLF section).\n#  5. New sub-bullet inserted after \"PSCAD version 5.0 used for results shown\"\n#     (before the E-TRAN bullet) explaining the Timed Breaker Logic switch,\n#     and the E-TRAN bullet's two runs (\"...PSCAD fil\" + \"e\") are\n#     consolidated into one run (no visible text change).\n\n$d = $word.ActiveDocument\n\n$LDQ = [char]0x201C   # \"\n$RDQ = [char]0x201D   # \"\n$LSQ = [char]0x2018   # '\n$RSQ = [char]0x2019   # '\n\nfunction Find-ParaByText($doc, $target) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t.Length -gt 0) {\n            $body = $t.Substring(0, $t.Length - 1)\n        } else {\n            $body = $t\n        }\n        if ($body -eq $target) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Set-ParaText($doc, $para, $newText) {\n    # Exclude the trailing paragraph mark from the replaced range, otherwise\n    # Range.Text = ... leaves the old tail behind instead of replacing it.\n    $r = $para.Range\n    $bodyRange = $doc.Range($r.Start, $r.End - 1)\n    $bodyRange.Text = $newText\n}\n\nfunction Insert-BulletAfter($doc, $afterPara, $newText) {\n    $r = $afterPara.Range\n    $r.Collapse(0)             # wdCollapseEnd\n    $r.InsertParagraphAfter()\n    $newPara = $afterPara.Next()\n    $newPara.Range.ListFormat.ListLevelNumber = 2   # matches ilvl=1 (0-based) of the sibling bullets\n    $newPara.Range.InsertBefore($newText)\n    return $newPara\n}\n\n# --- 1. Normalize the title paragraph to a single run -------------------\n$title = $d.Paragraphs(1)\nSet-ParaText $d $title \"READ ME: Enhanced IEEE 39 Bus System with Three IBRs\"\n\n# --- 2. Insert the eventType/disturbance.py bullet (PSSE section) -------\n$psseVersionPara = Find-ParaByText $d \"PSSE version 35 used for results shown \"\n$eventTypeText = \"By modifying the value of \" + $LDQ + \"eventType\" + $RDQ + \" in the \" + $LSQ + \"disturbance.py\" + $RSQ + \" file, one can simulate either a Busfault or GenTrip scenario\"\nInsert-BulletAfter $d $psseVersionPara $eventTypeText | Out-Null\n\n# --- 3. Normalize the EPCL (.p) bullet to a single run -------------------\n$epclPara = Find-ParaByText $d \"EPCL (.p) files with for running dynamic simulation with run commands\"\nSet-ParaText $d $epclPara \"EPCL (.p) files with for running dynamic simulation with run commands\"\n\n# --- 4. Insert the case_idx/Disturbance.p bullet (PSLF section) ---------\n$pslfVersionPara = Find-ParaByText $d \"PSLF version 32 used for results shown\"\n$caseIdxText = \"By modifying the value of \" + $LDQ + \"case_idx\" + $RDQ + \" in the \" + $LSQ + \"Disturbance.p\" + $RSQ + \" file, one can simulate either a BusFault or GenTrip scenario\"\nInsert-BulletAfter $d $pslfVersionPara $caseIdxText | Out-Null\n\n# --- 5. Insert the Timed Breaker Logic bullet (PSCAD section) + fix ETRAN bullet\n$pscadVersionPara = Find-ParaByText $d \"PSCAD version 5.0 used for results shown\"\n$breakerText = \"By modifying the time of breaker operation in \" + $LDQ + \"Timed Breaker Logic\" + $RDQ + \" in PSCAD, one can simulate either a BusFault or GenTrip scenario\"\nInsert-BulletAfter $d $pscadVersionPara $breakerText | Out-Null\n\n$etranPara = Find-ParaByText $d \"E-TRAN runtime library file for initializing ETRAN-based components in PSCAD file\"\nSet-ParaText $d $etranPara \"E-TRAN runtime library file for initializing ETRAN-based components in PSCAD file\"\n"}
